$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "24÷8=" "84÷2="
Replace-Text "94÷8=" "59÷6="
Replace-Text "17÷9=" "84÷5="
Replace-Text "48÷2=" "31÷3="
Replace-Text "78÷7=" "53÷6="
Replace-Text "64÷3=" "49÷2="
Replace-Text "34÷6=" "49÷2="
Replace-Text "63÷9=" "83÷6="
Replace-Text "13÷7=" "64÷7="
Replace-Text "38÷7=" "77÷3="
Replace-Text "72÷2=" "97÷9="
Replace-Text "15÷2=" "14÷8="
Replace-Text "21÷9=" "77÷9="
Replace-Text "18÷5=" "20÷4="
Replace-Text "35÷4=" "86÷8="
Replace-Text "78÷9=" "51÷6="
Replace-Text "29÷6=" "56÷3="
Replace-Text "20÷9=" "92÷4="
Replace-Text "39÷3=" "34÷7="
Replace-Text "65÷5=" "93÷7="
Replace-Text "24÷4=" "60÷3="
Replace-Text "37÷2=" "53÷5="
Replace-Text "72÷7=" "32÷2="
Replace-Text "33÷5=" "75÷3="
Replace-Text "21÷8=" "94÷2="
